$wb = $excel.ActiveWorkbook

# --- Sheet "展览" ---
$ws1 = $wb.Worksheets.Item("展览")

$ws1.Range("F3").Value = 55
$ws1.Range("F4").Value = 1470
$ws1.Range("F5").Value = 417
$ws1.Range("F6").Value = 1055
$ws1.Range("F7").Value = 10921
$ws1.Range("I9").Value = "//i2.hdslb.com/bfs/openplatform/202403/GWNvc78z1709275224442.jpeg"
$ws1.Range("F10").Value = 304
$ws1.Range("I10").Value = "//i0.hdslb.com/bfs/openplatform/202403/nIPoXWqO1709275656198.jpeg"
$ws1.Range("F12").Value = 738
$ws1.Range("F14").Value = 12666

# --- Sheet "全部类型" ---
$ws4 = $wb.Worksheets.Item("全部类型")

$ws4.Range("F4").Value = 55
$ws4.Range("F5").Value = 1470
$ws4.Range("F6").Value = 417
$ws4.Range("F7").Value = 1055
$ws4.Range("F8").Value = 10921
$ws4.Range("F9").Value = 10921
$ws4.Range("I11").Value = "//i2.hdslb.com/bfs/openplatform/202403/GWNvc78z1709275224442.jpeg"
$ws4.Range("F12").Value = 304
$ws4.Range("I12").Value = "//i0.hdslb.com/bfs/openplatform/202403/nIPoXWqO1709275656198.jpeg"
$ws4.Range("F14").Value = 738
$ws4.Range("F16").Value = 12666
